$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: mark the two "Expected HTTPS..." columns as <TBD> and highlight them in red
$ws.Range("F1").Value = "Expected HTTPS Code<TBD>"
$ws.Range("G1").Value = "Expected HTTPS Message/Body<TBD>"
$ws.Range("F1:G1").Font.Color = 255

# Point each "create customer" test row at its own dedicated JSON test-data file
# instead of the generic one they were all sharing before
$ws.Range("E2").Value = "customerCreationValid.json"
$ws.Range("E3").Value = "customerCreationInValidCustomerName.json"
$ws.Range("E4").Value = "customerCreationInValidEmail.json"
$ws.Range("E5").Value = "customerCreationInValidMobileNum.json"

# Clear the stray "Test Data" values that had leaked into the Expected-Message column
$ws.Range("G7").Clear()
$ws.Range("G8").ClearContents()

# Update the active selection
$ws.Range("G3").Select()
